$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cairo & Giza POIs")

# Replace empty / NaN "Entry cost (EGP)" cells (column G) with 0 for rows 211-260
$ws.Range("G211:G260").Value = 0

# Reflect the saved selection/view state on the sheet:
# scrolled back to the top (frozen header pane, topLeftCell = A2) with the
# final selection sitting just past the last data row.
$ws.Range("A2").Select()
$ws.Range("G261").Select()
